# Update the two-digit ÷ one-digit division problems/answers in the
# single table of this worksheet. Each "data" table row (1, 5, 9, 13, 17)
# holds 5 cells (columns 1-5); rows 2-4, 6-8, etc. are blank spacer rows.
#
# We address each cell directly via Tables.Item/Cell(row,col) and replace
# its text, instead of a global Find/Replace, because some new values are
# identical to other cells' old values (e.g. "79÷2=39, 1") - a global
# Find/Replace run sequentially would incorrectly re-match text that was
# just inserted by an earlier step.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $newText
}

# Data row 1 (table row 1)
Set-CellText $t 1 1 "92÷5=18, 2"
Set-CellText $t 1 2 "58÷4=14, 2"
Set-CellText $t 1 3 "79÷2=39, 1"
Set-CellText $t 1 4 "31÷7=4, 3"
Set-CellText $t 1 5 "21÷4=5, 1"

# Data row 2 (table row 5)
Set-CellText $t 5 1 "82÷5=16, 2"
Set-CellText $t 5 2 "98÷5=19, 3"
Set-CellText $t 5 3 "95÷6=15, 5"
Set-CellText $t 5 4 "62÷5=12, 2"
Set-CellText $t 5 5 "33÷9=3, 6"

# Data row 3 (table row 9)
Set-CellText $t 9 1 "80÷9=8, 8"
Set-CellText $t 9 2 "73÷6=12, 1"
Set-CellText $t 9 3 "59÷8=7, 3"
Set-CellText $t 9 4 "50÷2=25, 0"
Set-CellText $t 9 5 "54÷6=9, 0"

# Data row 4 (table row 13)
Set-CellText $t 13 1 "47÷5=9, 2"
Set-CellText $t 13 2 "75÷2=37, 1"
Set-CellText $t 13 3 "72÷6=12, 0"
Set-CellText $t 13 4 "63÷4=15, 3"
Set-CellText $t 13 5 "81÷5=16, 1"

# Data row 5 (table row 17)
Set-CellText $t 17 1 "33÷9=3, 6"
Set-CellText $t 17 2 "40÷8=5, 0"
Set-CellText $t 17 3 "96÷5=19, 1"
Set-CellText $t 17 4 "96÷7=13, 5"
Set-CellText $t 17 5 "76÷9=8, 4"
